$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Periodo Mora" detail value (E16: "2507" -> "2508")
$ws.Range("E16").Value = "2508"

# Update "Valor Mora" total (E11) and the matching detail value (F16):
# both go from 11388 to 56940
$ws.Range("E11").Value = 56940
$ws.Range("F16").Value = 56940
